$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, shifting the existing header row (now row 2)
# and all the data rows (now rows 3-32) down by one.
$ws.Rows.Item(1).Insert()

# New row 1: a sequence of numbers 0..11 across columns A..L, using the same
# bold/centered/bordered style that the header row previously used (style index 1,
# i.e. whatever style is currently applied to row 2, which used to be row 1).
$headerStyleRange = $ws.Range("A2:L2")
$ws.Cells.Item(1,1).Value = 0
$ws.Cells.Item(1,2).Value = 1
$ws.Cells.Item(1,3).Value = 2
$ws.Cells.Item(1,4).Value = 3
$ws.Cells.Item(1,5).Value = 4
$ws.Cells.Item(1,6).Value = 5
$ws.Cells.Item(1,7).Value = 6
$ws.Cells.Item(1,8).Value = 7
$ws.Cells.Item(1,9).Value = 8
$ws.Cells.Item(1,10).Value = 9
$ws.Cells.Item(1,11).Value = 10
$ws.Cells.Item(1,12).Value = 11

# Copy the style from the old header row (now shifted to row 2) onto new row 1
$headerStyleRange.Copy()
$newRowRange = $ws.Range("A1:L1")
$newRowRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Clear the text that used to be in I2, K2, L2 (now blank inline strings per the target)
$ws.Range("I2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
